# Update column F (dSF) values for the matching rows, per the diff:
# row -> new value
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -6
    3  = 2
    4  = 8
    5  = 7
    6  = -5
    9  = 3
    10 = -3
    11 = -3
    12 = -4
    13 = 4
    14 = 4
    15 = 3
    16 = -2
    17 = 1
    18 = 2
    19 = -2
    20 = 7
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
